# Update "建筑业企业总收入.xlsx" worksheet:
#  - Drop the 2000年-2009年 rows (old rows 2-10)
#  - Keep 2010年-2020年 (shift up into rows 2-11)
#  - Append a new 2021年 row (row 12)
#  - Dimension becomes A1:G12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 9 data rows (2000年 .. 2009年), shifting everything else up.
$ws.Range("A2:A10").EntireRow.Delete() | Out-Null

# Duplicate formatting of the last existing data row (now row 11, 2020年) onto
# the new row 12 so the new year label keeps the same cell style (s="1").
$ws.Range("A11").Copy($ws.Range("A12")) | Out-Null

# Populate the new 2021年 row.
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 8193.19623
$ws.Range("C12").Value = 241391.14125
$ws.Range("D12").Value = 262453.81074
$ws.Range("E12").Value = 239.8214
$ws.Range("F12").Value = 5441.92616
$ws.Range("G12").Value = 267895.7369
